# Refresh the recalculated ligand/receptor TPM-derived specificity metrics
# for the surviving (non-MuSCs-sending) rows, then drop the obsolete
# MuSCs-sending-cluster block (old rows 8-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.001904333333333333
$ws.Range("H2").Value = 0.005713
$ws.Range("I2").Value = 0.01334809965397277
$ws.Range("J2").Value = 0.01334809965397277
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 0.003860303299777778
$ws.Range("R2").Value = 0.034742729698
$ws.Range("S2").Value = 0.00008804786372517015
$ws.Range("T2").Value = 0.00008804786372517014

# Row 3
$ws.Range("G3").Value = 0.001904333333333333
$ws.Range("H3").Value = 0.005713
$ws.Range("I3").Value = 0.01334809965397277
$ws.Range("J3").Value = 0.01334809965397277
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 0.4883554387897778
$ws.Range("R3").Value = 4.395198949108
$ws.Range("S3").Value = 0.01113867221948164
$ws.Range("T3").Value = 0.01113867221948163

# Row 4
$ws.Range("G4").Value = 0.001904333333333333
$ws.Range("H4").Value = 0.005713
$ws.Range("I4").Value = 0.01334809965397277
$ws.Range("J4").Value = 0.01334809965397277
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 0.0930081459178889
$ws.Range("R4").Value = 0.8370733132610001
$ws.Range("S4").Value = 0.002121379570765966
$ws.Range("T4").Value = 0.002121379570765965

# Row 5
$ws.Range("I5").Value = 0.9866519003460271
$ws.Range("J5").Value = 0.9866519003460271
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 0.2853421599608889
$ws.Range("R5").Value = 2.568079439648
$ws.Range("S5").Value = 0.006508236701693444
$ws.Range("T5").Value = 0.006508236701693444

# Row 6
$ws.Range("I6").Value = 0.9866519003460271
$ws.Range("J6").Value = 0.9866519003460271
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("S6").Value = 0.8233375834448557
$ws.Range("T6").Value = 0.8233375834448556

# Row 7
$ws.Range("I7").Value = 0.9866519003460271
$ws.Range("J7").Value = 0.9866519003460271
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 6.874886035948446
$ws.Range("R7").Value = 61.873974323536
$ws.Range("S7").Value = 0.156806080199478
$ws.Range("T7").Value = 0.156806080199478

# Remove the obsolete MuSCs-sending-cluster rows (old rows 8-10)
$ws.Rows("8:10").Delete()

